$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999989451458875
$ws.Range("E2").Value = 0.9999989451458875

# Row 3
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = 0.4913720682015342
$ws.Range("E3").Value = 0.4913720682015342

# Row 4
$ws.Range("D4").Value = [double]"5.089818755836831E-27"
$ws.Range("E4").Value = [double]"5.089818755836831E-27"

# Row 5
$ws.Range("D5").Value = 0.1487322841097617
$ws.Range("E5").Value = 0.1487322841097617

# Row 6
$ws.Range("D6").Value = 0.9999964337274111
$ws.Range("E6").Value = 0.9999964337274111

# Row 7
$ws.Range("D7").Value = 0.002809850116856534
$ws.Range("E7").Value = 0.9971901498831435

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.4972219757779539
$ws.Range("E8").Value = 0.5027780242220461

# Row 10
$ws.Range("D10").Value = [double]"3.09653195388882E-05"
$ws.Range("E10").Value = 0.9999690346804612

# Row 11
$ws.Range("D11").Value = 0.9984680879944932
$ws.Range("E11").Value = 0.001531912005506841
$ws.Range("F11").Value = 4.410068035125732
